$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = '@'
$ws.Range("D2").Value = '27.964.63'
$ws.Range("D2").Style = 'Normal'
$ws.Range("E2").Value = '  -0.08%  '
$ws.Range("D3").NumberFormat = '@'
$ws.Range("D3").Value = '1.858.44'
$ws.Range("D3").Style = 'Normal'
$ws.Range("E3").Value = '  -1.33%  '
$ws.Range("D4").NumberFormat = '@'
$ws.Range("D4").Value = '1.004'
$ws.Range("D4").Style = 'Normal'
$ws.Range("E4").Value = '  +0.34%  '
$ws.Range("D5").NumberFormat = '@'
$ws.Range("D5").Value = '311.47'
$ws.Range("D5").Style = 'Normal'
$ws.Range("E5").Value = '  -0.44%  '
$ws.Range("E6").Value = '  +0.29%  '
$ws.Range("D7").NumberFormat = '@'
$ws.Range("D7").Value = '0.5115'
$ws.Range("D7").Style = 'Normal'
$ws.Range("E7").Value = '  +2.38%  '
$ws.Range("D8").NumberFormat = '@'
$ws.Range("D8").Value = '0.3804'
$ws.Range("D8").Style = 'Normal'
$ws.Range("E8").Value = '  -1.12%  '
$ws.Range("D9").NumberFormat = '@'
$ws.Range("D9").Value = '0.08286'
$ws.Range("D9").Style = 'Normal'
$ws.Range("E9").Value = '  -9.90%  '
$ws.Range("D10").NumberFormat = '@'
$ws.Range("D10").Value = '1.105'
$ws.Range("D10").Style = 'Normal'
$ws.Range("E10").Value = '  -1.37%  '
$ws.Range("D11").NumberFormat = '@'
$ws.Range("D11").Value = '41.24'
$ws.Range("D11").Style = 'Normal'
$ws.Range("E11").Value = '  -0.84%  '
$ws.Range("D12").NumberFormat = '@'
$ws.Range("D12").Value = '6.185'
$ws.Range("D12").Style = 'Normal'
$ws.Range("E12").Value = '  -2.31%  '
$ws.Range("D13").NumberFormat = '@'
$ws.Range("D13").Value = '20.41'
$ws.Range("D13").Style = 'Normal'
$ws.Range("E13").Value = '  -1.53%  '
$ws.Range("D14").NumberFormat = '@'
$ws.Range("D14").Value = '1.857.63'
$ws.Range("D14").Style = 'Normal'
$ws.Range("E14").Value = '  -0.96%  '
$ws.Range("D15").NumberFormat = '@'
$ws.Range("D15").Value = '7.170'
$ws.Range("D15").Style = 'Normal'
$ws.Range("E15").Value = '  -1.59%  '
$ws.Range("D16").NumberFormat = '@'
$ws.Range("D16").Value = '1.004'
$ws.Range("D16").Style = 'Normal'
$ws.Range("E16").Value = '  +0.39%  '
$ws.Range("D17").NumberFormat = '@'
$ws.Range("D17").Value = '0.00001093'
$ws.Range("D17").Style = 'Normal'
$ws.Range("E17").Value = '  -1.12%  '
$ws.Range("D18").NumberFormat = '@'
$ws.Range("D18").Value = '90.23'
$ws.Range("D18").Style = 'Normal'
$ws.Range("E18").Value = '  -1.17%  '
$ws.Range("D19").NumberFormat = '@'
$ws.Range("D19").Value = '0.06617'
$ws.Range("D19").Style = 'Normal'
$ws.Range("D20").NumberFormat = '@'
$ws.Range("D20").Value = '17.71'
$ws.Range("D20").Style = 'Normal'
$ws.Range("E20").Value = '  -1.48%  '
$ws.Range("E21").Value = '  +0.16%  '
$ws.Range("D22").NumberFormat = '@'
$ws.Range("D22").Value = '6.000'
$ws.Range("D22").Style = 'Normal'
$ws.Range("E22").Value = '  -2.71%  '
$ws.Range("D23").NumberFormat = '@'
$ws.Range("D23").Value = '28.005.05'
$ws.Range("D23").Style = 'Normal'
$ws.Range("E23").Value = '  -0.12%  '
$ws.Range("D24").NumberFormat = '@'
$ws.Range("D24").Value = '11.01'
$ws.Range("D24").Style = 'Normal'
$ws.Range("E24").Value = '  -3.29%  '
$ws.Range("D25").NumberFormat = '@'
$ws.Range("D25").Value = '2.258'
$ws.Range("D25").Style = 'Normal'
$ws.Range("E25").Value = '  -2.13%  '
$ws.Range("D26").NumberFormat = '@'
$ws.Range("D26").Value = '2.562'
$ws.Range("D26").Style = 'Normal'
$ws.Range("E26").Value = '  +0.84%  '
$ws.Range("D27").NumberFormat = '@'
$ws.Range("D27").Value = '2.071.21'
$ws.Range("D27").Style = 'Normal'
$ws.Range("E27").Value = '  -1.04%  '
$ws.Range("D28").NumberFormat = '@'
$ws.Range("D28").Value = '157.12'
$ws.Range("D28").Style = 'Normal'
$ws.Range("E28").Value = '  -0.52%  '
$ws.Range("D29").NumberFormat = '@'
$ws.Range("D29").Value = '20.37'
$ws.Range("D29").Style = 'Normal'
$ws.Range("E29").Value = '  -1.89%  '
$ws.Range("D30").NumberFormat = '@'
$ws.Range("D30").Value = '124.51'
$ws.Range("D30").Style = 'Normal'
$ws.Range("E30").Value = '  -1.60%  '
$ws.Range("D31").NumberFormat = '@'
$ws.Range("D31").Value = '0.1057'
$ws.Range("D31").Style = 'Normal'
$ws.Range("E31").Value = '  +0.25%  '
$ws.Range("D32").NumberFormat = '@'
$ws.Range("D32").Value = '1.038'
$ws.Range("D32").Style = 'Normal'
$ws.Range("E32").Value = '  -2.71%  '
$ws.Range("D33").NumberFormat = '@'
$ws.Range("D33").Value = '5.575'
$ws.Range("D33").Style = 'Normal'
$ws.Range("E33").Value = '  -0.16%  '
$ws.Range("D34").NumberFormat = '@'
$ws.Range("D34").Value = '3.603'
$ws.Range("D34").Style = 'Normal'
$ws.Range("E34").Value = '  +0.48%  '
$ws.Range("D35").NumberFormat = '@'
$ws.Range("D35").Value = '9.632'
$ws.Range("D35").Style = 'Normal'
$ws.Range("E35").Value = '  +2.84%  '
$ws.Range("D36").NumberFormat = '@'
$ws.Range("D36").Value = '0.06508'
$ws.Range("D36").Style = 'Normal'
$ws.Range("E36").Value = '  -0.93%  '
$ws.Range("D37").NumberFormat = '@'
$ws.Range("D37").Value = '0.02405'
$ws.Range("D37").Style = 'Normal'
$ws.Range("E37").Value = '  +0.26%  '
$ws.Range("D38").NumberFormat = '@'
$ws.Range("D38").Value = '0.2151'
$ws.Range("D38").Style = 'Normal'
$ws.Range("E38").Value = '  -1.40%  '
$ws.Range("D39").NumberFormat = '@'
$ws.Range("D39").Value = '1.205'
$ws.Range("D39").Style = 'Normal'
$ws.Range("E39").Value = '  +0.05%  '
$ws.Range("D40").NumberFormat = '@'
$ws.Range("D40").Value = '0.6386'
$ws.Range("D40").Style = 'Normal'
$ws.Range("E40").Value = '  -0.38%  '
$ws.Range("D41").NumberFormat = '@'
$ws.Range("D41").Value = '1.225'
$ws.Range("D41").Style = 'Normal'
$ws.Range("E41").Value = '  -4.58%  '
$ws.Range("D42").NumberFormat = '@'
$ws.Range("D42").Value = '11.23'
$ws.Range("D42").Style = 'Normal'
$ws.Range("E42").Value = '  -3.21%  '
$ws.Range("D43").NumberFormat = '@'
$ws.Range("D43").Value = '4.861'
$ws.Range("D43").Style = 'Normal'
$ws.Range("E43").Value = '  -1.50%  '
$ws.Range("D44").NumberFormat = '@'
$ws.Range("D44").Value = '0.6039'
$ws.Range("D44").Style = 'Normal'
$ws.Range("E44").Value = '  +0.34%  '
$ws.Range("D45").NumberFormat = '@'
$ws.Range("D45").Value = '13.06'
$ws.Range("D45").Style = 'Normal'
$ws.Range("E45").Value = '  -1.74%  '
$ws.Range("D46").NumberFormat = '@'
$ws.Range("D46").Value = '1.284'
$ws.Range("D46").Style = 'Normal'
$ws.Range("E46").Value = '  -0.66%  '
$ws.Range("E47").Value = '  -0.45%  '
$ws.Range("D48").NumberFormat = '@'
$ws.Range("D48").Value = '1.968'
$ws.Range("D48").Style = 'Normal'
$ws.Range("E48").Value = '  -1.00%  '
$ws.Range("D49").NumberFormat = '@'
$ws.Range("D49").Value = '1.205'
$ws.Range("D49").Style = 'Normal'
$ws.Range("E49").Value = '  +0.44%  '
$ws.Range("D50").NumberFormat = '@'
$ws.Range("D50").Value = '120.33'
$ws.Range("D50").Style = 'Normal'
$ws.Range("E50").Value = '  -0.80%  '
$ws.Range("D51").NumberFormat = '@'
$ws.Range("D51").Value = '79.55'
$ws.Range("D51").Style = 'Normal'
$ws.Range("E51").Value = '  +0.93%  '
